$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 16 detail values (LINDAO ZUÑIGA BRYAN JOSE / CHASI PASTO ANGEL NOLBERTO)
$ws1.Range("H16").Value = 772.64
$ws1.Range("I16").Value = 739.5
$ws1.Range("M16").Value = 612.86
$ws1.Range("N16").Value = 2042.06

# Row 58 summary counters ("X de 56")
$ws1.Range("H58").Value = "1 de 56"
$ws1.Range("I58").Value = "2 de 56"
$ws1.Range("M58").Value = "7 de 56"
$ws1.Range("N58").Value = "2 de 56"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F16").Value = 4167.06
$ws2.Range("F58").Value = 29573.5

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 7 - INODOROS
$ws3.Range("D7").Value = 772.64
$ws3.Range("E7").Value = 1327.36
$ws3.Range("F7").Value = 0.3679238095238095

# Row 8 - LAVABOS
$ws3.Range("D8").Value = 825.9
$ws3.Range("E8").Value = -75.89999999999998
$ws3.Range("F8").Value = 1.1012

# Row 16 - PORCELANATO
$ws3.Range("D16").Value = 4371.46
$ws3.Range("E16").Value = 34045.71
$ws3.Range("F16").Value = 0.1137892249741457

# Row 17 - PUERTAS DE SEGURIDAD
$ws3.Range("D17").Value = 2268.66
$ws3.Range("E17").Value = -1926.66
$ws3.Range("F17").Value = 6.633508771929824

# Row 19 - TOTAL
$ws3.Range("D19").Value = 29573.5
$ws3.Range("E19").Value = 25836.20560036207
$ws3.Range("F19").Value = 0.5337241856741927
